# Updated cryptos list on Mon Dec  4 10:45:43 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for rows 2-51.
# Every D/E cell in this sheet is stored as text (not a number), so numeric
# looking Price values (e.g. "234.53") are written with a leading "'" to
# stop Excel from auto-coercing them into Number cells, then the cell
# Style is reset back to "Normal" so no stray quote-prefix formatting is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.968.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.19%  "
$ws.Range("D3").Value = "'2.262.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'234.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.40%  "
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("D7").Value = "'63.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.411"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.84%  "
$ws.Range("D10").Value = "'60.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("D11").Value = "'0.0899"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").Value = "'2.597.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.35%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'22.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.72%  "
$ws.Range("D16").Value = "'0.823"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "'5.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.41%  "
$ws.Range("D18").Value = "'2.263.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("D19").Value = "'41.778.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.79%  "
$ws.Range("D20").Value = "'74.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("D21").Value = "'0.0₃0935"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.13%  "
$ws.Range("D22").Value = "'6.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "
$ws.Range("D23").Value = "'252.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +10.16%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").Value = "'2.35"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +7.87%  "
$ws.Range("D28").Value = "'9.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D29").Value = "'170.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").Value = "'20.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("D31").Value = "'1.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("E32").Value = "  +7.80%  "
$ws.Range("E33").Value = "  +1.07%  "
$ws.Range("D34").Value = "'5.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.42%  "
$ws.Range("E35").Value = "  +3.78%  "
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("D37").Value = "'6.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").Value = "'2.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("D40").Value = "'0.000264"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +57.74%  "
$ws.Range("D41").Value = "'5.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.27%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("E43").Value = "  +6.05%  "
$ws.Range("D44").Value = "'8.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.86%  "
$ws.Range("D45").Value = "'17.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "'0.0987"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.75%  "
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("D49").Value = "'1.504.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "
$ws.Range("E50").Value = "  +1.44%  "
$ws.Range("D51").Value = "'2.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.20%  "
